# Apply the "Updated Bar Graph Percent file" edit to the Summary slide (slide 8):
#  1. Move the Title placeholder to its new position.
#  2. Remove the Slide Number placeholder from this slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- 1. Reposition the title shape ---
# Target offset in EMU: x=676798, y=4240203.
# Shape.Left / Shape.Top are expressed in points (1 pt = 12700 EMU); add a
# half-EMU bias before the conversion so the value lands safely inside the
# target EMU's rounding interval once PowerPoint stores it back as EMU.
$title = $s.Shapes.Item("Title 10")
$title.Left = (676798 + 0.5) / 12700
$title.Top  = (4240203 + 0.5) / 12700

# --- 2. Remove the Slide Number placeholder shape from this slide ---
# Turning off slide-number visibility for the slide drops the placeholder
# shape from the slide (same as unchecking it in Insert > Header & Footer).
$s.HeadersFooters.SlideNumber.Visible = $false
